# LOT2041.docx edit: insert manual line breaks into two paragraphs.
#
# 1) "Programa resumido" paragraph: fix the missing space between
#    "Formas de" and "condução" by splitting the run at that point and
#    inserting a manual line break (<w:br/>) there.
#
# 2) "Bibliografia" paragraph: split the single run into three runs,
#    each one reference, separated by manual line breaks (<w:br/>).

$d = $word.ActiveDocument

# --- Edit 1: Programa resumido -------------------------------------------
$old1 = "Características do material biológico; Cinética de processos " + `
        "fermentativos, Formas decondução dos processos fermentativos, " + `
        "esterilização em bioprocessos."
$new1 = "Características do material biológico; Cinética de processos " + `
        "fermentativos, Formas de^lcondução dos processos fermentativos, " + `
        "esterilização em bioprocessos."

$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, `
                                   $false, $true, 1, $false, $new1, 2)
Write-Host "Edit 1 applied: $found1"

# --- Edit 2: Bibliografia -------------------------------------------------
$dash = [char]8211   # "–" EN DASH

$ref1 = "1.Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia " + `
        "Industrial $dash Engenharia Bioquímica, vol. 2, São Paulo: Edgard " + `
        "Blücher, 2001. "
$ref2 = "2. Borzani, W.; Schmidell, W.; Lima, U. A.; Aquarone, E. Biotecnologia " + `
        "Industrial. Fundamentos Vol. 1. São Paulo: Ed. Edgard Blücher, 2001. "
$ref3 = "3. Pauline M. Doran ed. Bioprocess Engineering Principles " + `
        "(Second Edition), Elsevier Ltd. 2013."

$old2 = $ref1 + $ref2 + $ref3
$new2 = $ref1 + "^l" + $ref2 + "^l" + $ref3

$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, `
                                   $false, $true, 1, $false, $new2, 2)
Write-Host "Edit 2 applied: $found2"
